# New daily price observation is inserted for Feria Lagunitas de Puerto Montt - Pomelo.
# The new record belongs right before the existing row 33 (date-sorted data), so every
# row from 33 downward shifts down by one (old row 33 -> new row 34, ..., old row 150
# -> new row 151), and the freshly inserted row 33 carries the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 33..150 down to 34..151, leaving a blank row 33 (keeps column formatting,
# e.g. the date style on column D, in sync with the rest of the table).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C33").Value = "Los Lagos"
$ws.Range("D33").Value = 44481
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100102
$ws.Range("H33").Value = "Cítricos"
$ws.Range("I33").Value = 100102006
$ws.Range("J33").Value = "Pomelo"
$ws.Range("K33").Value = "Start Ruby"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 160
$ws.Range("N33").Value = 11000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 11500
$ws.Range("Q33").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R33").Value = "Región de O'Higgins"
$ws.Range("S33").Value = 821
$ws.Range("T33").Value = 14
